$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: demographic_clinical
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("demographic_clinical")
$ws1.Range("A2").Value = "/mnt/munin/Morey/Lab/Delin/Projects/IBMMA/Data/ENIGMA-PGC_master_v1.1_TR&fID_v1.xlsx"
$ws1.Columns.Item(1).ColumnWidth = 76.5
$ws1.Range("A5").Select()

# ---------------------------------------------------------------------------
# Sheet 2: data_path
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("data_path")
# Drop the old "atlas_conn" row and keep only the falff_reho path (updated)
$ws2.Rows.Item(2).Delete()
$ws2.Range("A2").Value = "/mnt/munin/Morey/Lab/Delin/Projects/IBMMA/Data/falff_reho"
$ws2.Columns.Item(1).ColumnWidth = 50.166666666666664
$ws2.Range("A8").Select()

# ---------------------------------------------------------------------------
# Sheet 3: data_pattern
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("data_pattern")

# New headers
$ws3.Range("E1").Value = "MASK1"
$ws3.Range("F1").Value = "ROI1"
$ws3.Range("G1").Value = "ROI2"
$ws3.Range("H1").Value = "MYROI"
$ws3.Range("I1").Value = "MEASURE"
$ws3.Range("J1").Value = "EXCLUDED"

# Row 2 (fALFF_alff) gains a MASK1 path and moves its EXCLUDED flag to column J
$ws3.Range("E2").Value = "/mnt/munin/Morey/Lab/Delin/Projects/IBMMA/Data/brain_mask.nii"
$ws3.Range("G2").ClearContents()
$ws3.Range("J2").Value = 1

# Row 3: replace the old corrMatrix entry with the new reHo entry
$ws3.Range("A3").Value = "reHo"
$ws3.Range("B3").Value = "_feature-reHo_reho.nii.gz"
$ws3.Range("C3").Value = "_feature-reHo_reho.json"
$ws3.Range("D3").Value = "_feature-reHo_mask.nii.gz"
$ws3.Range("E3").Value = "/mnt/munin/Morey/Lab/Delin/Projects/IBMMA/Data/brain_mask.nii"
$ws3.Range("G3").ClearContents()
$ws3.Range("J3").Value = 0

# Column widths
$ws3.Columns.Item(1).ColumnWidth = 10.5
$ws3.Columns.Item(2).ColumnWidth = 21
$ws3.Columns.Item(3).ColumnWidth = 20.666666666666668
$ws3.Columns.Item(4).ColumnWidth = 24
$ws3.Columns.Item(5).ColumnWidth = 53.83333333333333
$ws3.Columns.Item(6).ColumnWidth = 6.833333333333334
$ws3.Columns.Item(7).ColumnWidth = 5.5
$ws3.Columns.Item(8).ColumnWidth = 6.5
$ws3.Columns.Item(9).ColumnWidth = 9.333333333333332

$ws3.Activate()
$ws3.Range("E9").Select()

# ---------------------------------------------------------------------------
# Sheet 4: predictors -- unchanged content, nothing to do
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("predictors")
$ws4.Range("A2").Select()

# ---------------------------------------------------------------------------
# Sheet 5: models -- add Model_03
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("models")
$ws5.Range("A4").Value = "Model_03"
$ws5.Range("B4").Value = "lmer(Yvar ~ GROUP * SEX + AGE +  (1|SITE))"
$ws5.Range("C4").Value = 0
$ws5.Range("B11").Select()
